# CAN1-40 Fixed error on FIFO requirements by Antonio
# - Duplication on FIFO_08 and FIFO_11 requirements

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: paragraph "FIFO_00" - the sentence was split across a stray
# _GoBack bookmark ("...that can h" | "ave a value ..."). Re-join the text so
# it reads naturally and drop the bookmark that split it.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("that can have a value", $true, $false, $false, $false, $false, $true, 1, $false, "that can have a value", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: paragraph "FIFO_08" is missing the "unless o_full is 1" clause
# that every other similar requirement (FIFO_07, FIFO_09, FIFO_11) has.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("when i_w_en is 1.", $true, $false, $false, $false, $false, $true, 1, $false, "when i_w_en is 1, unless o_full is 1.", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 3: the requirement "FIFO_11" paragraph was duplicated verbatim.
# Remove the duplicate paragraph entirely.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("FIFO_11: The register the write pointer")) {
        $p.Range.Delete()
        break
    }
}

# The following paragraph was "FIFO_12" and must be renumbered to "FIFO_11"
# now that the duplicate has been removed.
$d.Content.Find.Execute("FIFO_12:", $true, $false, $false, $false, $false, $true, 1, $false, "FIFO_11:", 2) | Out-Null

# Re-add the _GoBack bookmark at the very end of the document (its new home
# after the edit, per the last-saved cursor position).
$endRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$d.Bookmarks.Add("_GoBack", $endRange) | Out-Null
